$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1: remove the extra-tall row height (back to default) ---
$ws.Rows.Item(1).AutoFit()

# --- Fill in "Id" (A), "Start time" (B), "Completion time" (C), and
#     "Email" (D, anonymized) for the newly-tested respondents (rows 2-7) ---
$ws.Range("A2").Value = 2
$ws.Range("B2").Value = 45067.989583333336
$ws.Range("C2").Value = 45067.990972222222
$ws.Range("D2").Value = "anonymous"

$ws.Range("A3").Value = 3
$ws.Range("B3").Value = 45069.840277777781
$ws.Range("C3").Value = 45069.842361111114
$ws.Range("D3").Value = "anonymous"

$ws.Range("A4").Value = 4
$ws.Range("B4").Value = 45071.589583333334
$ws.Range("C4").Value = 45071.592361111114
$ws.Range("D4").Value = "anonymous"

$ws.Range("A5").Value = 5
$ws.Range("B5").Value = 45071.646527777775
$ws.Range("C5").Value = 45071.648611111108
$ws.Range("D5").Value = "anonymous"

$ws.Range("A6").Value = 6
$ws.Range("B6").Value = 45071.915277777778
$ws.Range("C6").Value = 45071.919444444444
$ws.Range("D6").Value = "anonymous"

$ws.Range("A7").Value = 7
$ws.Range("B7").Value = 45072.915277777778
$ws.Range("C7").Value = 45072.919444444444
$ws.Range("D7").Value = "anonymous"

# --- Update the view: scroll back to column A and move the active selection ---
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("D14").Select()
